$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.194.14"
$ws.Range("E2").Value = "'  -3.08%  "
$ws.Range("D3").Value = "'1.649.15"
$ws.Range("E3").Value = "'  -3.37%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "'  +0.22%  "
$ws.Range("D5").Value = "'309.13"
$ws.Range("E5").Value = "'  -2.17%  "
$ws.Range("E6").Value = "'  +0.13%  "
$ws.Range("D7").Value = "'0.3898"
$ws.Range("E7").Value = "'  -1.41%  "
$ws.Range("D8").Value = "'0.3887"
$ws.Range("E8").Value = "'  -3.48%  "
$ws.Range("D9").Value = "'1.003"
$ws.Range("E9").Value = "'  +0.17%  "
$ws.Range("D10").Value = "'1.372"
$ws.Range("E10").Value = "'  -7.48%  "
$ws.Range("D11").Value = "'49.45"
$ws.Range("E11").Value = "'  -6.15%  "
$ws.Range("D12").Value = "'0.08500"
$ws.Range("E12").Value = "'  -3.62%  "
$ws.Range("E13").Value = "'  -6.59%  "
$ws.Range("E14").Value = "'  -3.85%  "
$ws.Range("E15").Value = "'  -4.60%  "
$ws.Range("D16").Value = "'7.540"
$ws.Range("E16").Value = "'  -5.71%  "
$ws.Range("D17").Value = "'1.650.70"
$ws.Range("E17").Value = "'  -3.82%  "
$ws.Range("D18").Value = "'94.99"
$ws.Range("E18").Value = "'  -1.32%  "
$ws.Range("B19").Value = "'Avalanche"
$ws.Range("C19").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'21.26"
$ws.Range("E19").Value = "'  +3.22%  "
$ws.Range("B20").Value = "'TRON"
$ws.Range("C20").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D20").Value = "'0.06919"
$ws.Range("E20").Value = "'  -3.67%  "
$ws.Range("D21").Value = "'6.980"
$ws.Range("E21").Value = "'  -5.11%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "'  +0.19%  "
$ws.Range("D23").Value = "'13.87"
$ws.Range("E23").Value = "'  -4.07%  "
$ws.Range("D24").Value = "'24.193.61"
$ws.Range("E24").Value = "'  -3.05%  "
$ws.Range("D25").Value = "'2.369"
$ws.Range("E25").Value = "'  +0.68%  "
$ws.Range("D26").Value = "'2.765"
$ws.Range("E27").Value = "'  -4.35%  "
$ws.Range("D28").Value = "'158.59"
$ws.Range("E28").Value = "'  -1.94%  "
$ws.Range("D29").Value = "'8.594"
$ws.Range("E29").Value = "'  +1.99%  "
$ws.Range("D30").Value = "'143.52"
$ws.Range("E30").Value = "'  -4.65%  "
$ws.Range("D31").Value = "'5.362"
$ws.Range("E31").Value = "'  -13.99%  "
$ws.Range("D32").Value = "'2.415"
$ws.Range("E32").Value = "'  -7.58%  "
$ws.Range("D33").Value = "'1.832.05"
$ws.Range("E33").Value = "'  -4.10%  "
$ws.Range("D34").Value = "'6.996"
$ws.Range("E34").Value = "'  -2.88%  "
$ws.Range("D35").Value = "'0.08120"
$ws.Range("E35").Value = "'  -5.28%  "
$ws.Range("D36").Value = "'0.9975"
$ws.Range("E36").Value = "'  -4.66%  "
$ws.Range("D37").Value = "'0.02940"
$ws.Range("E37").Value = "'  -6.04%  "
$ws.Range("D38").Value = "'0.2722"
$ws.Range("E38").Value = "'  -4.72%  "
$ws.Range("D39").Value = "'0.09334"
$ws.Range("E39").Value = "'  -2.16%  "
$ws.Range("D40").Value = "'1.482"
$ws.Range("E40").Value = "'  -0.57%  "
$ws.Range("D41").Value = "'10.04"
$ws.Range("E41").Value = "'  -7.38%  "
$ws.Range("D42").Value = "'0.7685"
$ws.Range("E42").Value = "'  -6.96%  "
$ws.Range("D43").Value = "'13.24"
$ws.Range("E43").Value = "'  -5.52%  "
$ws.Range("D44").Value = "'16.17"
$ws.Range("E44").Value = "'  -6.86%  "
$ws.Range("D45").Value = "'2.514"
$ws.Range("E45").Value = "'  -6.63%  "
$ws.Range("D46").Value = "'0.6930"
$ws.Range("E46").Value = "'  -6.21%  "
$ws.Range("D47").Value = "'4.100"
$ws.Range("E47").Value = "'  -3.55%  "
$ws.Range("E48").Value = "'  +0.04%  "
$ws.Range("D49").Value = "'0.08463"
$ws.Range("E49").Value = "'  -3.41%  "
$ws.Range("D50").Value = "'1.278"
$ws.Range("E50").Value = "'  -9.84%  "
$ws.Range("D51").Value = "'134.92"
$ws.Range("E51").Value = "'  -3.07%  "
